# Apply "Actualización automática de tasas-transfi.xlsx" changes

$wb = $excel.ActiveWorkbook

# --- 1. Update the conversion summary text on sheet "Hoja1", cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.51 = 17823.91 pesos`n✅ 17823.91 pesos = 4.47 = 956.18 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- 2. Update tasas rates on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 221.949
$wsTasas.Range("O10").Value = 3956
$wsTasas.Range("N12").Value = 3989.12
$wsTasas.Range("O12").Value = 214
